# The host pre-declares $excel / $wb / $app, but in this runtime $wb starts
# out un-bound, so re-resolve it from $excel explicitly before using it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary row (row 49): case-sensitive accuracy header + values ---
$ws.Range("A49").Value = "case sensitive"
$ws.Range("B49").Value = "accuracy on train"
$ws.Range("C49").Value = "accuracy on test"

# --- Column B got wider to fit the new "accuracy on train" header ---
# Range/Columns.ColumnWidth is in Excel "characters" units and gets
# pixel-quantized by the host, so this lands on the closest representable
# width to the authored 16.3877551020408 (≈ 15.5544 chars -> 16.3333 width).
$ws.Columns.Item(2).ColumnWidth = 15.554421768707465

# --- View state: scroll the window down and select the new row ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.TabRatio = 0.142
$ws.Range("A49").Select() | Out-Null

Write-Host "Added row 49 (case sensitive / accuracy on train / accuracy on test)"
